$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.851.15"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").Value = "'2.093.89"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'234.08"
$ws.Range("E5").Value = "  -0.65%  "
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").Value = "'58.72"
$ws.Range("E7").Value = "  -0.42%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  +0.42%  "
$ws.Range("E10").Value = "  -0.90%  "
$ws.Range("E11").Value = "  +3.03%  "
$ws.Range("D12").Value = "'15.24"
$ws.Range("E12").Value = "  +1.95%  "
$ws.Range("D13").Value = "'2.401.97"
$ws.Range("E13").Value = "  +0.26%  "
$ws.Range("D14").Value = "'21.37"
$ws.Range("E14").Value = "  +0.34%  "
$ws.Range("D15").Value = "'0.782"
$ws.Range("E15").Value = "  +0.70%  "
$ws.Range("D16").Value = "'5.38"
$ws.Range("E16").Value = "  +0.87%  "
$ws.Range("D17").Value = "'2.075.15"
$ws.Range("E17").Value = "  -0.61%  "
$ws.Range("D18").Value = "'37.802.64"
$ws.Range("E18").Value = "  +0.23%  "
$ws.Range("E19").Value = "  -0.51%  "
$ws.Range("E20").Value = "  -0.36%  "
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("D22").Value = "'230.60"
$ws.Range("E22").Value = "  +0.66%  "
$ws.Range("D23").Value = "'0.998"
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("D24").Value = "'2.40"
$ws.Range("E24").Value = "  -0.43%  "
$ws.Range("E25").Value = "  -1.28%  "
$ws.Range("D26").Value = "'9.84"
$ws.Range("E26").Value = "  +8.55%  "
$ws.Range("D27").Value = "'171.30"
$ws.Range("E27").Value = "  +0.78%  "
$ws.Range("E28").Value = "  -2.83%  "
$ws.Range("D29").Value = "'19.56"
$ws.Range("E29").Value = "  -0.44%  "
$ws.Range("D30").Value = "'1.41"
$ws.Range("E30").Value = "  -0.50%  "
$ws.Range("E31").Value = "  +0.18%  "
$ws.Range("D32").Value = "'4.73"
$ws.Range("E32").Value = "  +0.08%  "
$ws.Range("D33").Value = "'0.0635"
$ws.Range("E33").Value = "  -0.23%  "
$ws.Range("E34").Value = "  -0.58%  "
$ws.Range("E35").Value = "  -0.22%  "
$ws.Range("E36").Value = "  -1.00%  "
$ws.Range("D37").Value = "'3.33"
$ws.Range("E37").Value = "  -3.53%  "
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("D39").Value = "'5.40"
$ws.Range("E39").Value = "  -0.51%  "
$ws.Range("D40").Value = "'0.0236"
$ws.Range("E40").Value = "  +9.13%  "
$ws.Range("D41").Value = "'101.92"
$ws.Range("E41").Value = "  +2.83%  "
$ws.Range("E42").Value = "  -0.71%  "
$ws.Range("E43").Value = "  +0.84%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").Value = "'16.79"
$ws.Range("E44").Value = "  +4.61%  "
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").Value = "'1.18"
$ws.Range("E45").Value = "  +0.85%  "
$ws.Range("D46").Value = "'1.451.71"
$ws.Range("E46").Value = "  -0.83%  "
$ws.Range("E47").Value = "  -5.65%  "
$ws.Range("E48").Value = "  -0.50%  "
$ws.Range("D49").Value = "'7.28"
$ws.Range("E49").Value = "  -2.62%  "
$ws.Range("E50").Value = "  -2.03%  "
$ws.Range("D51").Value = "'2.285.02"
$ws.Range("E51").Value = "  +0.21%  "
